$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-12 22:18:28"
$ws.Range("E3").Value = "2026-02-12 22:18:31"
$ws.Range("E4").Value = "2026-02-12 22:18:33"
$ws.Range("J4").Value = "999.8 hPa"
$ws.Range("O4").Value = "16.1 °C"
$ws.Range("E5").Value = "2026-02-12 22:18:36"
$ws.Range("E6").Value = "2026-02-12 22:18:38"
$ws.Range("J6").Value = "999.6 hPa"
$ws.Range("O6").Value = "15.7 °C"
$ws.Range("E7").Value = "2026-02-12 22:18:41"
$ws.Range("J7").Value = "1002.3 hPa"
$ws.Range("E8").Value = "2026-02-12 22:18:43"
$ws.Range("J8").Value = "1001.7 hPa"
$ws.Range("E9").Value = "2026-02-12 22:18:46"
$ws.Range("E10").Value = "2026-02-12 22:18:48"
$ws.Range("O10").Value = "14.7 °C"
$ws.Range("E11").Value = "2026-02-12 22:18:51"
$ws.Range("H11").Value = "'48%"
$ws.Range("O11").Value = "9.2 °C"
$ws.Range("E12").Value = "2026-02-12 22:18:53"
$ws.Range("O12").Value = "12.5 °C"
$ws.Range("E13").Value = "2026-02-12 22:18:56"
$ws.Range("J13").Value = "1002.3 hPa"
$ws.Range("O13").Value = "7.5 °C"
$ws.Range("E14").Value = "2026-02-12 22:18:58"
$ws.Range("N14").Value = "12.6 °C 21:42 TU"
$ws.Range("O14").Value = "16.9 °C"
$ws.Range("E15").Value = "2026-02-12 22:19:01"
$ws.Range("E16").Value = "2026-02-12 22:19:03"
$ws.Range("E17").Value = "2026-02-12 22:19:06"
$ws.Range("H17").Value = "'75%"
$ws.Range("E18").Value = "2026-02-12 22:19:08"
$ws.Range("H18").Value = "'38%"
$ws.Range("J18").Value = "1000.0 hPa"
$ws.Range("N18").Value = "9.7 °C 21:58 TU"
$ws.Range("O18").Value = "16.5 °C"
$ws.Range("E19").Value = "2026-02-12 22:19:11"
$ws.Range("O19").Value = "7.9 °C"
$ws.Range("E20").Value = "2026-02-12 22:19:14"
$ws.Range("E21").Value = "2026-02-12 22:19:16"
$ws.Range("J21").Value = "1002.8 hPa"
$ws.Range("O21").Value = "9.0 °C"
$ws.Range("E22").Value = "2026-02-12 22:19:19"
$ws.Range("E23").Value = "2026-02-12 22:19:21"
$ws.Range("E24").Value = "2026-02-12 22:19:24"
$ws.Range("O24").Value = "11.5 °C"
$ws.Range("E25").Value = "2026-02-12 22:19:27"
$ws.Range("E26").Value = "2026-02-12 22:19:29"
$ws.Range("J26").Value = "999.4 hPa"
$ws.Range("N26").Value = "2.1 °C 21:48 TU"
$ws.Range("O26").Value = "5.8 °C"
$ws.Range("E27").Value = "2026-02-12 22:19:31"
$ws.Range("E28").Value = "2026-02-12 22:19:34"
$ws.Range("H28").Value = "'40%"
$ws.Range("J28").Value = "999.5 hPa"
$ws.Range("O28").Value = "13.7 °C"
$ws.Range("E29").Value = "2026-02-12 22:19:36"
$ws.Range("H29").Value = "'61%"
$ws.Range("N29").Value = "5.9 °C 21:30 TU"
$ws.Range("O29").Value = "14.0 °C"
$ws.Range("E30").Value = "2026-02-12 22:19:39"
$ws.Range("E31").Value = "2026-02-12 22:19:41"
$ws.Range("J31").Value = "999.3 hPa"
$ws.Range("E32").Value = "2026-02-12 22:19:44"
$ws.Range("E33").Value = "2026-02-12 22:19:46"
$ws.Range("J33").Value = "1001.9 hPa"
$ws.Range("O33").Value = "6.5 °C"
$ws.Range("E34").Value = "2026-02-12 22:19:49"
$ws.Range("O34").Value = "0.3 °C"
$ws.Range("E35").Value = "2026-02-12 22:19:52"
$ws.Range("J35").Value = "1008.2 hPa"
$ws.Range("E36").Value = "2026-02-12 22:19:54"
$ws.Range("H36").Value = "'58%"
$ws.Range("J36").Value = "1000.1 hPa"
$ws.Range("E37").Value = "2026-02-12 22:19:57"
$ws.Range("H37").Value = "'49%"
$ws.Range("J37").Value = "1000.9 hPa"
$ws.Range("N37").Value = "3.0 °C 21:58 TU"
$ws.Range("O37").Value = "9.7 °C"
$ws.Range("E38").Value = "2026-02-12 22:19:59"
$ws.Range("E39").Value = "2026-02-12 22:20:02"
$ws.Range("E40").Value = "2026-02-12 22:20:04"
$ws.Range("H40").Value = "'56%"
$ws.Range("J40").Value = "1003.6 hPa"
$ws.Range("N40").Value = "4.6 °C 21:59 TU"
$ws.Range("O40").Value = "9.3 °C"
$ws.Range("E41").Value = "2026-02-12 22:20:07"
$ws.Range("E42").Value = "2026-02-12 22:20:09"
$ws.Range("H42").Value = "'62%"
$ws.Range("N42").Value = "7.5 °C 21:37 TU"
$ws.Range("O42").Value = "13.8 °C"
$ws.Range("E43").Value = "2026-02-12 22:20:12"
$ws.Range("E44").Value = "2026-02-12 22:20:14"
$ws.Range("E45").Value = "2026-02-12 22:20:17"
$ws.Range("J45").Value = "1005.4 hPa"
$ws.Range("N45").Value = "2.6 °C 21:53 TU"
$ws.Range("O45").Value = "6.8 °C"
$ws.Range("E46").Value = "2026-02-12 22:20:19"
$ws.Range("J46").Value = "1007.6 hPa"
